$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

$ws.Range("O2:O12").Value = "nan"
